# Auto-generated cell value updates for Sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 17700
$ws.Range("J7").Value = 10500
$ws.Range("L7").Value = 10500
$ws.Range("N7").Value = -10724
$ws.Range("H14").Value = 17700
$ws.Range("J14").Value = 10500
$ws.Range("L14").Value = 10500
$ws.Range("N14").Value = -10882
$ws.Range("H86").Value = 100002480
$ws.Range("I86").Value = 100002480
$ws.Range("K86").Value = 100002480
$ws.Range("M86").Value = -100001357
$ws.Range("H89").Value = 100002480
$ws.Range("I89").Value = 100002480
$ws.Range("K89").Value = 500012400
$ws.Range("M89").Value = -500006784
$ws.Range("H107").Value = 1080.3684
$ws.Range("I107").Value = 1053.7142
$ws.Range("J107").Value = 1155
$ws.Range("K107").Value = 1053.7142
$ws.Range("L107").Value = 1155
$ws.Range("M107").Value = 866.2858000000001
$ws.Range("N107").Value = -4995
$ws.Range("H113").Value = 3594.8
$ws.Range("J113").Value = 3594.8
$ws.Range("L113").Value = 3594.8
$ws.Range("N113").Value = -10102.8
$ws.Range("H141").Value = 2199.5
$ws.Range("I141").Value = 2199.5
$ws.Range("K141").Value = 6598.5
$ws.Range("M141").Value = -1418.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19670.256
$ws.Range("I32").Value = 19145.6
$ws.Range("K32").Value = 19145.6
$ws.Range("M32").Value = -18858.6
$ws.Range("H45").Value = 3146.6667
$ws.Range("J45").Value = 4060
$ws.Range("L45").Value = 4060
$ws.Range("N45").Value = -4814
$ws.Range("H61").Value = 1817.25
$ws.Range("I61").Value = 1817.25
$ws.Range("K61").Value = 1817.25
$ws.Range("M61").Value = -1605.25
$ws.Range("H63").Value = 928
$ws.Range("I63").Value = 950
$ws.Range("J63").Value = 906
$ws.Range("K63").Value = 950
$ws.Range("L63").Value = 906
$ws.Range("M63").Value = -264
$ws.Range("N63").Value = -2278
$ws.Range("H66").Value = 928
$ws.Range("I66").Value = 950
$ws.Range("J66").Value = 906
$ws.Range("K66").Value = 4750
$ws.Range("L66").Value = 4530
$ws.Range("M66").Value = -1318
$ws.Range("N66").Value = -11394
$ws.Range("H74").Value = 1366.9131
$ws.Range("I74").Value = 1173.5294
$ws.Range("K74").Value = 1173.5294
$ws.Range("M74").Value = -299.5293999999999
$ws.Range("H77").Value = 1366.9131
$ws.Range("I77").Value = 1173.5294
$ws.Range("K77").Value = 5867.646999999999
$ws.Range("M77").Value = -1499.646999999999
$ws.Range("H132").Value = 2214.0715
$ws.Range("I132").Value = 2214.0715
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6642.2145
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4112.2145
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 1817.25
$ws.Range("I136").Value = 1817.25
$ws.Range("K136").Value = 5451.75
$ws.Range("M136").Value = -2901.75

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8391.385
$ws.Range("I20").Value = 7669.1
$ws.Range("K20").Value = 7669.1
$ws.Range("M20").Value = -7422.1
$ws.Range("H29").Value = 750
$ws.Range("I29").Value = 750
$ws.Range("K29").Value = 750
$ws.Range("M29").Value = -461
$ws.Range("H105").Value = 2616.4666
$ws.Range("I105").Value = 2590.4167
$ws.Range("J105").Value = 2720.6667
$ws.Range("K105").Value = 2590.4167
$ws.Range("L105").Value = 2720.6667
$ws.Range("M105").Value = -843.4167000000002
$ws.Range("N105").Value = -6214.6667
$ws.Range("H107").Value = 3143.5417
$ws.Range("I107").Value = 3250.7856
$ws.Range("J107").Value = 2993.4
$ws.Range("K107").Value = 3250.7856
$ws.Range("L107").Value = 2993.4
$ws.Range("M107").Value = -1330.7856
$ws.Range("N107").Value = -6833.4
$ws.Range("H134").Value = 2909.5625
$ws.Range("I134").Value = 2228.4443
$ws.Range("J134").Value = 3785.2856
$ws.Range("K134").Value = 6685.3329
$ws.Range("L134").Value = 11355.8568
$ws.Range("M134").Value = -4150.3329
$ws.Range("N134").Value = -16425.8568

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1037.4
$ws.Range("I16").Value = 874.8889
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 874.8889
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -587.8889
$ws.Range("N16").Value = -3074
$ws.Range("H31").Value = 23812352
$ws.Range("I31").Value = 25643458
$ws.Range("K31").Value = 25643458
$ws.Range("M31").Value = -25643163
$ws.Range("H34").Value = 23812352
$ws.Range("I34").Value = 25643458
$ws.Range("K34").Value = 25643458
$ws.Range("M34").Value = -25643256
$ws.Range("H113").Value = 1037.4
$ws.Range("I113").Value = 874.8889
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 874.8889
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 1295.1111
$ws.Range("N113").Value = -6840
$ws.Range("H125").Value = 60000
$ws.Range("J125").Value = 60000
$ws.Range("L125").Value = 60000
$ws.Range("N125").Value = -64920
$ws.Range("H132").Value = 121229140
$ws.Range("I132").Value = 222232590
$ws.Range("K132").Value = 666697770
$ws.Range("M132").Value = -666695240

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 464.7143
$ws.Range("I15").Value = 170
$ws.Range("J15").Value = 685.75
$ws.Range("K15").Value = 510
$ws.Range("L15").Value = 2057.25
$ws.Range("M15").Value = -370
$ws.Range("N15").Value = -2337.25
$ws.Range("H68").Value = 1413.2222
$ws.Range("J68").Value = 1575
$ws.Range("L68").Value = 4725
$ws.Range("N68").Value = -6347
$ws.Range("H69").Value = 5097.3
$ws.Range("I69").Value = 4395.2
$ws.Range("K69").Value = 13185.6
$ws.Range("M69").Value = -12374.6
$ws.Range("H71").Value = 1413.2222
$ws.Range("J71").Value = 1575
$ws.Range("L71").Value = 14175
$ws.Range("N71").Value = -22287
$ws.Range("H72").Value = 5097.3
$ws.Range("I72").Value = 4395.2
$ws.Range("K72").Value = 39556.8
$ws.Range("M72").Value = -35500.8
$ws.Range("H74").Value = 6499.5
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 7999
$ws.Range("K74").Value = 15000
$ws.Range("L74").Value = 23997
$ws.Range("M74").Value = -13939
$ws.Range("N74").Value = -26119
$ws.Range("H75").Value = 1906.3334
$ws.Range("I75").Value = 359.5
$ws.Range("J75").Value = 5000
$ws.Range("K75").Value = 1078.5
$ws.Range("L75").Value = 15000
$ws.Range("M75").Value = -80.5
$ws.Range("N75").Value = -16996
$ws.Range("H76").Value = 6124.75
$ws.Range("I76").Value = 1999
$ws.Range("K76").Value = 5997
$ws.Range("M76").Value = -5614
$ws.Range("H77").Value = 6499.5
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 7999
$ws.Range("K77").Value = 45000
$ws.Range("L77").Value = 71991
$ws.Range("M77").Value = -39696
$ws.Range("N77").Value = -82599
$ws.Range("H78").Value = 1906.3334
$ws.Range("I78").Value = 359.5
$ws.Range("J78").Value = 5000
$ws.Range("K78").Value = 3235.5
$ws.Range("L78").Value = 45000
$ws.Range("M78").Value = 1756.5
$ws.Range("N78").Value = -54984
$ws.Range("H79").Value = 6124.75
$ws.Range("I79").Value = 1999
$ws.Range("K79").Value = 5997
$ws.Range("M79").Value = -4671
$ws.Range("H87").Value = 18989.334
$ws.Range("J87").Value = 18989.334
$ws.Range("L87").Value = 56968.00199999999
$ws.Range("N87").Value = -59464.00199999999
$ws.Range("H90").Value = 18989.334
$ws.Range("J90").Value = 18989.334
$ws.Range("L90").Value = 170904.006
$ws.Range("N90").Value = -183384.006

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 10602200
$ws.Range("I14").Value = 13251750
$ws.Range("K14").Value = 13251750
$ws.Range("M14").Value = -13251582
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 459979.28
$ws.Range("I122").Value = 2503305
$ws.Range("K122").Value = 7509915
$ws.Range("M122").Value = -7507465
$ws.Range("H132").Value = 127350.19
$ws.Range("I132").Value = 155099.16
$ws.Range("J132").Value = 7104.6665
$ws.Range("K132").Value = 465297.48
$ws.Range("L132").Value = 21313.9995
$ws.Range("M132").Value = -462767.48
$ws.Range("N132").Value = -26373.9995
$ws.Range("H141").Value = 32451.4
$ws.Range("J141").Value = 32451.4
$ws.Range("L141").Value = 32451.4
$ws.Range("N141").Value = -42811.4

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4111.703
$ws.Range("I132").Value = 3390.509
$ws.Range("J132").Value = 5213.528
$ws.Range("K132").Value = 10171.527
$ws.Range("L132").Value = 15640.584
$ws.Range("M132").Value = -7641.527
$ws.Range("N132").Value = -20700.584

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 944.1667
$ws.Range("I4").Value = 933
$ws.Range("K4").Value = 933
$ws.Range("M4").Value = -820
$ws.Range("H132").Value = 1566.527
$ws.Range("I132").Value = 1347.3846
$ws.Range("K132").Value = 4042.1538
$ws.Range("M132").Value = -1512.1538
$ws.Range("H140").Value = 95273
$ws.Range("J140").Value = 95273
$ws.Range("L140").Value = 95273
$ws.Range("N140").Value = -105633
$ws.Range("H141").Value = 94000
$ws.Range("J141").Value = 94000
$ws.Range("L141").Value = 94000
$ws.Range("N141").Value = -104360

